# Fix data files names
# Replace ASCII-art console output (tidyverse startup message) with the
# actual Unicode glyphs R/RStudio renders (box-drawing dashes and
# check/cross/info symbols) in the "SourceCode" verbatim block.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "-- Attaching core tidyverse packages ------------------------ tidyverse 2.0.0 --" "── Attaching core tidyverse packages ──────────────────────── tidyverse 2.0.0 ──"
Replace-Text "v dplyr     1.1.4     v readr     2.1.6" "✔ dplyr     1.1.4     ✔ readr     2.1.6"
Replace-Text "v forcats   1.0.1     v stringr   1.6.0" "✔ forcats   1.0.1     ✔ stringr   1.6.0"
Replace-Text "v ggplot2   4.0.1     v tibble    3.3.1" "✔ ggplot2   4.0.1     ✔ tibble    3.3.1"
Replace-Text "v lubridate 1.9.4     v tidyr     1.3.2" "✔ lubridate 1.9.4     ✔ tidyr     1.3.2"
Replace-Text "v purrr     1.2.1     " "✔ purrr     1.2.1     "
Replace-Text "-- Conflicts ------------------------------------------ tidyverse_conflicts() --" "── Conflicts ────────────────────────────────────────── tidyverse_conflicts() ──"
Replace-Text "x dplyr::filter() masks stats::filter()" "✖ dplyr::filter() masks stats::filter()"
Replace-Text "x dplyr::lag()    masks stats::lag()" "✖ dplyr::lag()    masks stats::lag()"
Replace-Text "i Use the conflicted package (<http://conflicted.r-lib.org/>) to force all conflicts to become errors" "ℹ Use the conflicted package (<http://conflicted.r-lib.org/>) to force all conflicts to become errors"
